$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-11-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-08 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("37×83=3071", $true, $false, $false, $false, $false, $true, 1, $false, "71×19=1349", 2) | Out-Null
$d.Content.Find.Execute("73×30=2190", $true, $false, $false, $false, $false, $true, 1, $false, "75×74=5550", 2) | Out-Null
$d.Content.Find.Execute("96×74=7104", $true, $false, $false, $false, $false, $true, 1, $false, "69×52=3588", 2) | Out-Null
$d.Content.Find.Execute("14×15=210", $true, $false, $false, $false, $false, $true, 1, $false, "92×40=3680", 2) | Out-Null
$d.Content.Find.Execute("56×19=1064", $true, $false, $false, $false, $false, $true, 1, $false, "35×24=840", 2) | Out-Null
$d.Content.Find.Execute("84×41=3444", $true, $false, $false, $false, $false, $true, 1, $false, "92×67=6164", 2) | Out-Null
$d.Content.Find.Execute("48×20=960", $true, $false, $false, $false, $false, $true, 1, $false, "97×31=3007", 2) | Out-Null
$d.Content.Find.Execute("85×32=2720", $true, $false, $false, $false, $false, $true, 1, $false, "96×84=8064", 2) | Out-Null
$d.Content.Find.Execute("30×32=960", $true, $false, $false, $false, $false, $true, 1, $false, "89×20=1780", 2) | Out-Null
$d.Content.Find.Execute("48×53=2544", $true, $false, $false, $false, $false, $true, 1, $false, "50×54=2700", 2) | Out-Null
$d.Content.Find.Execute("71×18=1278", $true, $false, $false, $false, $false, $true, 1, $false, "67×12=804", 2) | Out-Null
$d.Content.Find.Execute("37×64=2368", $true, $false, $false, $false, $false, $true, 1, $false, "83×85=7055", 2) | Out-Null
$d.Content.Find.Execute("36×99=3564", $true, $false, $false, $false, $false, $true, 1, $false, "22×17=374", 2) | Out-Null
$d.Content.Find.Execute("95×44=4180", $true, $false, $false, $false, $false, $true, 1, $false, "25×58=1450", 2) | Out-Null
$d.Content.Find.Execute("74×23=1702", $true, $false, $false, $false, $false, $true, 1, $false, "68×56=3808", 2) | Out-Null
$d.Content.Find.Execute("66×84=5544", $true, $false, $false, $false, $false, $true, 1, $false, "64×39=2496", 2) | Out-Null
$d.Content.Find.Execute("36×64=2304", $true, $false, $false, $false, $false, $true, 1, $false, "30×76=2280", 2) | Out-Null
$d.Content.Find.Execute("60×11=660", $true, $false, $false, $false, $false, $true, 1, $false, "15×38=570", 2) | Out-Null
$d.Content.Find.Execute("42×68=2856", $true, $false, $false, $false, $false, $true, 1, $false, "25×71=1775", 2) | Out-Null
$d.Content.Find.Execute("55×12=660", $true, $false, $false, $false, $false, $true, 1, $false, "47×84=3948", 2) | Out-Null
$d.Content.Find.Execute("93×79=7347", $true, $false, $false, $false, $false, $true, 1, $false, "44×71=3124", 2) | Out-Null
$d.Content.Find.Execute("79×48=3792", $true, $false, $false, $false, $false, $true, 1, $false, "35×88=3080", 2) | Out-Null
$d.Content.Find.Execute("91×12=1092", $true, $false, $false, $false, $false, $true, 1, $false, "24×61=1464", 2) | Out-Null
$d.Content.Find.Execute("43×12=516", $true, $false, $false, $false, $false, $true, 1, $false, "30×54=1620", 2) | Out-Null
$d.Content.Find.Execute("57×23=1311", $true, $false, $false, $false, $false, $true, 1, $false, "67×32=2144", 2) | Out-Null
